$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 3499.1428
$ws.Range("I18").Value = 3499.1428
$ws.Range("K18").Value = 3499.1428
$ws.Range("M18").Value = -3215.1428
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("N46").ClearContents()
$ws.Range("H58").Value = 732.5
$ws.Range("J58").Value = 2700
$ws.Range("L58").Value = 8100
$ws.Range("N58").Value = -8400
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("N60").ClearContents()
$ws.Range("H114").Value = 87500
$ws.Range("J114").Value = 85000
$ws.Range("L114").Value = 85000
$ws.Range("N114").Value = -93678
$ws.Range("H132").Value = 1325.9166
$ws.Range("I132").Value = 1268.6666
$ws.Range("K132").Value = 3805.9998
$ws.Range("M132").Value = -1275.9998
$ws.Range("H135").Value = 659.25
$ws.Range("I135").Value = 566.0526
$ws.Range("K135").Value = 5094.4734
$ws.Range("M135").Value = -2559.4734
$ws.Range("H137").Value = 2427.375
$ws.Range("I137").Value = 2172
$ws.Range("K137").Value = 6516
$ws.Range("M137").Value = -3966
$ws.Range("H138").Value = 3573.0317
$ws.Range("I138").Value = 2971.2144
$ws.Range("J138").Value = 3744.9795
$ws.Range("K138").Value = 8913.643199999999
$ws.Range("L138").Value = 11234.9385
$ws.Range("M138").Value = -3773.643199999999
$ws.Range("N138").Value = -21514.9385
$ws.Range("H141").Value = 1489.3334
$ws.Range("I141").Value = 1734.75
$ws.Range("K141").Value = 5204.25
$ws.Range("M141").Value = -24.25

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4066.9546
$ws.Range("I61").Value = 1344.9231
$ws.Range("K61").Value = 1344.9231
$ws.Range("M61").Value = -1132.9231
$ws.Range("I63").Value = 2000
$ws.Range("K63").Value = 2000
$ws.Range("M63").Value = -1314
$ws.Range("I66").Value = 2000
$ws.Range("K66").Value = 10000
$ws.Range("M66").Value = -6568
$ws.Range("H74").Value = 1901.4
$ws.Range("I74").Value = 1374.4193
$ws.Range("K74").Value = 1374.4193
$ws.Range("M74").Value = -500.4193
$ws.Range("H77").Value = 1901.4
$ws.Range("I77").Value = 1374.4193
$ws.Range("K77").Value = 6872.0965
$ws.Range("M77").Value = -2504.0965
$ws.Range("H92").Value = 60723.5
$ws.Range("J92").Value = 60723.5
$ws.Range("L92").Value = 60723.5
$ws.Range("N92").Value = -65715.5
$ws.Range("H102").Value = 1272.6666
$ws.Range("I102").Value = 1387.2
$ws.Range("J102").Value = 700
$ws.Range("K102").Value = 1387.2
$ws.Range("L102").Value = 700
$ws.Range("M102").Value = 234.8
$ws.Range("N102").Value = -3944
$ws.Range("H132").Value = 2052.6667
$ws.Range("I132").Value = 1989.2727
$ws.Range("K132").Value = 5967.8181
$ws.Range("M132").Value = -3437.8181
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H136").Value = 4066.9546
$ws.Range("I136").Value = 1344.9231
$ws.Range("K136").Value = 4034.7693
$ws.Range("M136").Value = -1484.7693

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 684.7273
$ws.Range("J22").Value = 833.3333
$ws.Range("L22").Value = 833.3333
$ws.Range("N22").Value = -1179.3333

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2309.3572
$ws.Range("I58").Value = 2194.3333
$ws.Range("K58").Value = 2194.3333
$ws.Range("M58").Value = -1991.3333
$ws.Range("H95").Value = 11042.75
$ws.Range("J95").Value = 11042.75
$ws.Range("L95").Value = 11042.75
$ws.Range("N95").Value = -16534.75
$ws.Range("H99").Value = 2857.5625
$ws.Range("I99").Value = 1938.3334
$ws.Range("K99").Value = 1938.3334
$ws.Range("M99").Value = -440.3334
$ws.Range("H105").Value = 3451.4
$ws.Range("J105").Value = 4282.8
$ws.Range("L105").Value = 4282.8
$ws.Range("N105").Value = -7776.8
$ws.Range("H107").Value = 1097.2222
$ws.Range("I107").Value = 699.3333
$ws.Range("K107").Value = 699.3333
$ws.Range("M107").Value = 1220.6667
$ws.Range("H126").Value = 2857.5625
$ws.Range("I126").Value = 1938.3334
$ws.Range("K126").Value = 5815.0002
$ws.Range("M126").Value = -3345.0002
$ws.Range("H134").Value = 1900.25
$ws.Range("I134").Value = 1925.85
$ws.Range("J134").Value = 1772.25
$ws.Range("K134").Value = 5777.549999999999
$ws.Range("L134").Value = 5316.75
$ws.Range("M134").Value = -3242.549999999999
$ws.Range("N134").Value = -10386.75
$ws.Range("H136").Value = 2309.3572
$ws.Range("I136").Value = 2194.3333
$ws.Range("K136").Value = 6582.999899999999
$ws.Range("M136").Value = -4032.999899999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 17142944
$ws.Range("I4").Value = 18333434
$ws.Range("K4").Value = 55000302
$ws.Range("M4").Value = -55000190
$ws.Range("H8").Value = 237.55556
$ws.Range("I8").Value = 237.55556
$ws.Range("K8").Value = 712.66668
$ws.Range("M8").Value = -573.66668
$ws.Range("H14").Value = 1177.9231
$ws.Range("I14").Value = 1177.9231
$ws.Range("K14").Value = 3533.7693
$ws.Range("M14").Value = -3360.7693
$ws.Range("H68").Value = 20840500
$ws.Range("J68").Value = 20840500
$ws.Range("L68").Value = 62521500
$ws.Range("N68").Value = -62523122
$ws.Range("H71").Value = 20840500
$ws.Range("J71").Value = 20840500
$ws.Range("L71").Value = 187564500
$ws.Range("N71").Value = -187572612
$ws.Range("H86").Value = 8333
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 8333
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 24999
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -27371
$ws.Range("H89").Value = 8333
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 8333
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 74997
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -86853
$ws.Range("H132").Value = 3820.8333
$ws.Range("I132").Value = 5032.5
$ws.Range("J132").Value = 2609.1667
$ws.Range("K132").Value = 45292.5
$ws.Range("L132").Value = 23482.5003
$ws.Range("M132").Value = -42762.5
$ws.Range("N132").Value = -28542.5003

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1146.0952
$ws.Range("I102").Value = 786.41174
$ws.Range("K102").Value = 786.41174
$ws.Range("M102").Value = 835.58826
$ws.Range("H123").Value = 7500
$ws.Range("J123").Value = 7500
$ws.Range("L123").Value = 7500
$ws.Range("N123").Value = -12400
$ws.Range("H132").Value = 2490.3333
$ws.Range("I132").Value = 2490.3333
$ws.Range("K132").Value = 7470.999899999999
$ws.Range("M132").Value = -4940.999899999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6838.5
$ws.Range("J7").Value = 9265.799999999999
$ws.Range("L7").Value = 9265.799999999999
$ws.Range("N7").Value = -9489.799999999999
$ws.Range("H22").Value = 1529.6
$ws.Range("I22").Value = 948.5
$ws.Range("J22").Value = 1917
$ws.Range("K22").Value = 948.5
$ws.Range("L22").Value = 1917
$ws.Range("M22").Value = -653.5
$ws.Range("N22").Value = -2507
$ws.Range("H27").Value = 1529.6
$ws.Range("I27").Value = 948.5
$ws.Range("J27").Value = 1917
$ws.Range("K27").Value = 948.5
$ws.Range("L27").Value = 1917
$ws.Range("M27").Value = -841.5
$ws.Range("N27").Value = -2131
$ws.Range("H80").Value = 60116
$ws.Range("I80").Value = 60116
$ws.Range("K80").Value = 60116
$ws.Range("M80").Value = -58993
$ws.Range("H83").Value = 60116
$ws.Range("I83").Value = 60116
$ws.Range("K83").Value = 180348
$ws.Range("M83").Value = -174732
$ws.Range("H94").Value = 49998
$ws.Range("J94").Value = 49998
$ws.Range("L94").Value = 49998
$ws.Range("N94").Value = -51350
$ws.Range("H126").Value = 6838.5
$ws.Range("J126").Value = 9265.799999999999
$ws.Range("L126").Value = 27797.4
$ws.Range("N126").Value = -32737.4
$ws.Range("H132").Value = 4502.5
$ws.Range("I132").Value = 4000
$ws.Range("K132").Value = 12000
$ws.Range("M132").Value = -9470
$ws.Range("H136").Value = 3985.7
$ws.Range("I136").Value = 3701.2666
$ws.Range("J136").Value = 4839
$ws.Range("K136").Value = 11103.7998
$ws.Range("L136").Value = 14517
$ws.Range("M136").Value = -8553.799800000001
$ws.Range("N136").Value = -19617

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 11926.8
$ws.Range("I55").Value = 4797
$ws.Range("K55").Value = 4797
$ws.Range("M55").Value = -4520
$ws.Range("H132").Value = 1553.4286
$ws.Range("I132").Value = 1553.4286
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4660.2858
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -2130.2858
$ws.Range("N132").ClearContents()
